$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G: hour changes from 22 to 23 for every data row (2-51)
$gRng = $ws.Range("G2:G51")
$gRng.NumberFormat = "@"
$gRng.Value = '23'
$gRng.Style = "Normal"

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.78'
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '25.22'
$ws.Range("D3").Style = "Normal"

# Row 4
$ws.Range("B4").Value = 'HuobiToken'
$ws.Range("C4").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.042'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '3HuobiTokenHT'

# Row 5
$ws.Range("B5").Value = 'Cronos'
$ws.Range("C5").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05595'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '4CronosCRO'

# Row 6
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.554'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '5KuCoinTokenKCS'

# Row 7
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.017'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '6GateTokenGT'

# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8140'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '7MXTokenMX'

# Row 9
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8345'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '8FTXTokenFTT'

# Row 10
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0005958'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '9OneONE'

# Row 11
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1335'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '10WazirXWRX'

# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06959'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'

# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.02827'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '12BitrueCoinBTR'

# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09404'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '13BitMartTokenBMX'

# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001514'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '14BitForexTokenBF'

# Row 16
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006089'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '15TigerCashTCH'

# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.500'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '16LEOLEO'

# Row 18
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.092'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '17BTSETokenBTSE'

# Row 19
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3187'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '18BitpandaEcosystemTokenBEST'

# Row 20
$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03250'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '19LiechtensteinCryptoassetsExchangeLCX'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.742'
$ws.Range("D22").Style = "Normal"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04689'
$ws.Range("D23").Style = "Normal"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009697'
$ws.Range("D27").Style = "Normal"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001939'
$ws.Range("D28").Style = "Normal"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03672'
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006238'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '40KickTokenKICKBestin24h'

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1058'
$ws.Range("D42").Style = "Normal"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002599'
$ws.Range("D43").Style = "Normal"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008206'
$ws.Range("D44").Style = "Normal"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005295'
$ws.Range("D45").Style = "Normal"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1799'
$ws.Range("D47").Style = "Normal"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002015'
$ws.Range("D48").Style = "Normal"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002099'
$ws.Range("D49").Style = "Normal"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'
$ws.Range("D50").Style = "Normal"
